$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated statistical description values per diff

# Row 2: co2
$ws.Range("C2").Value = 556.0118747503398
$ws.Range("D2").Value = 136.9457962890611
$ws.Range("F2").Value = 451
$ws.Range("G2").Value = 513
$ws.Range("H2").Value = 628

# Row 3: humidity
$ws.Range("C3").Value = 39.52368681086599
$ws.Range("D3").Value = 5.606686868718929
$ws.Range("F3").Value = 36.2
$ws.Range("G3").Value = 39.38
$ws.Range("H3").Value = 42.73

# Row 4: pm25
$ws.Range("C4").Value = 1.871659445233377
$ws.Range("D4").Value = 2.547662974226376
$ws.Range("F4").Value = 0.58
$ws.Range("G4").Value = 1.14
$ws.Range("H4").Value = 2.27

# Row 5: pressure
$ws.Range("C5").Value = 323.6522252399128
$ws.Range("D5").Value = 11.24365883105295
$ws.Range("F5").Value = 316.56
$ws.Range("G5").Value = 325.26
$ws.Range("H5").Value = 332.32

# Row 6: temperature
$ws.Range("C6").Value = 20.85940448499746
$ws.Range("D6").Value = 2.514335804458083
$ws.Range("F6").Value = 19.45
$ws.Range("G6").Value = 20.79
$ws.Range("H6").Value = 22.29

# Row 7: rssi
$ws.Range("C7").Value = -76.34696996049395
$ws.Range("D7").Value = 23.05853417398344
$ws.Range("H7").Value = -58

# Row 8: snr
$ws.Range("C8").Value = 7.530281690140845
$ws.Range("D8").Value = 6.973085129280163

# Row 9: SF
$ws.Range("C9").Value = 9.321963796632716
$ws.Range("D9").Value = 1.685559257197546

# Row 10: frequency
$ws.Range("C10").Value = 867.8300477245315
$ws.Range("D10").Value = 0.4614970413315863

# Row 11: toa
$ws.Range("C11").Value = 0.5556898985686461
$ws.Range("D11").Value = 0.5889505257655291

# Row 12: distance
$ws.Range("C12").Value = 22.7405931832473
$ws.Range("D12").Value = 12.2918382459891

# Row 13: c_walls
$ws.Range("C13").Value = 0.6741066223760586
$ws.Range("D13").Value = 0.7505689471373156

# Row 14: w_walls
$ws.Range("C14").Value = 1.826884373608353
$ws.Range("D14").Value = 1.663848617717044

# Row 15: exp_pl
$ws.Range("C15").Value = 93.74696996049376
$ws.Range("D15").Value = 23.05853417398344
$ws.Range("F15").Value = 75.40000000000001

# Row 16: n_power
$ws.Range("C16").Value = -85.51812118107686
$ws.Range("D16").Value = 20.69690774443287
$ws.Range("F16").Value = -102.4668316388797
$ws.Range("H16").Value = -68.8707776445072

# Row 17: esp
$ws.Range("C17").Value = -77.98783949093604
$ws.Range("D17").Value = 25.51063393418753
$ws.Range("F17").Value = -93.53779541063678
$ws.Range("G17").Value = -72.57382219273629
$ws.Range("H17").Value = -58.22214159641585
